$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet1 (Train Results) updates: rows 2-44, columns A-J
$ws1.Cells.Item(2,1).Value = 0 ; $ws1.Cells.Item(2,2).Value = 40 ; $ws1.Cells.Item(2,3).Value = 4 ; $ws1.Cells.Item(2,4).Value = 0 ; $ws1.Cells.Item(2,5).Value = 28 ; $ws1.Cells.Item(2,6).Value = 4 ; $ws1.Cells.Item(2,7).Value = 4 ; $ws1.Cells.Item(2,8).Value = 20 ; $ws1.Cells.Item(2,9).Value = 2.31 ; $ws1.Cells.Item(2,10).Value = 2.342777013778687
$ws1.Cells.Item(3,1).Value = 4 ; $ws1.Cells.Item(3,2).Value = 0 ; $ws1.Cells.Item(3,3).Value = 8 ; $ws1.Cells.Item(3,4).Value = 24 ; $ws1.Cells.Item(3,5).Value = 4 ; $ws1.Cells.Item(3,6).Value = 16 ; $ws1.Cells.Item(3,7).Value = 44 ; $ws1.Cells.Item(3,8).Value = 0 ; $ws1.Cells.Item(3,9).Value = 2.98 ; $ws1.Cells.Item(3,10).Value = 2.970501184463501
$ws1.Cells.Item(4,1).Value = 4 ; $ws1.Cells.Item(4,2).Value = 16 ; $ws1.Cells.Item(4,3).Value = 0 ; $ws1.Cells.Item(4,4).Value = 20 ; $ws1.Cells.Item(4,5).Value = 4 ; $ws1.Cells.Item(4,6).Value = 12 ; $ws1.Cells.Item(4,7).Value = 44 ; $ws1.Cells.Item(4,8).Value = 0 ; $ws1.Cells.Item(4,9).Value = 3.22 ; $ws1.Cells.Item(4,10).Value = 3.102986812591553
$ws1.Cells.Item(5,1).Value = 0 ; $ws1.Cells.Item(5,2).Value = 24 ; $ws1.Cells.Item(5,3).Value = 4 ; $ws1.Cells.Item(5,4).Value = 0 ; $ws1.Cells.Item(5,5).Value = 8 ; $ws1.Cells.Item(5,6).Value = 4 ; $ws1.Cells.Item(5,7).Value = 36 ; $ws1.Cells.Item(5,8).Value = 24 ; $ws1.Cells.Item(5,9).Value = 2.45 ; $ws1.Cells.Item(5,10).Value = 2.350285053253174
$ws1.Cells.Item(6,1).Value = 4 ; $ws1.Cells.Item(6,2).Value = 0 ; $ws1.Cells.Item(6,3).Value = 8 ; $ws1.Cells.Item(6,4).Value = 4 ; $ws1.Cells.Item(6,5).Value = 24 ; $ws1.Cells.Item(6,6).Value = 12 ; $ws1.Cells.Item(6,7).Value = 44 ; $ws1.Cells.Item(6,8).Value = 4 ; $ws1.Cells.Item(6,9).Value = 3.17 ; $ws1.Cells.Item(6,10).Value = 2.999998092651367
$ws1.Cells.Item(7,1).Value = 0 ; $ws1.Cells.Item(7,2).Value = 16 ; $ws1.Cells.Item(7,3).Value = 8 ; $ws1.Cells.Item(7,4).Value = 4 ; $ws1.Cells.Item(7,5).Value = 16 ; $ws1.Cells.Item(7,6).Value = 20 ; $ws1.Cells.Item(7,7).Value = 36 ; $ws1.Cells.Item(7,8).Value = 0 ; $ws1.Cells.Item(7,9).Value = 3.07 ; $ws1.Cells.Item(7,10).Value = 2.911064624786377
$ws1.Cells.Item(8,1).Value = 4 ; $ws1.Cells.Item(8,2).Value = 20 ; $ws1.Cells.Item(8,3).Value = 4 ; $ws1.Cells.Item(8,4).Value = 4 ; $ws1.Cells.Item(8,5).Value = 16 ; $ws1.Cells.Item(8,6).Value = 0 ; $ws1.Cells.Item(8,7).Value = 52 ; $ws1.Cells.Item(8,8).Value = 0 ; $ws1.Cells.Item(8,9).Value = 3.35 ; $ws1.Cells.Item(8,10).Value = 3.116973876953125
$ws1.Cells.Item(9,1).Value = 4 ; $ws1.Cells.Item(9,2).Value = 0 ; $ws1.Cells.Item(9,3).Value = 12 ; $ws1.Cells.Item(9,4).Value = 4 ; $ws1.Cells.Item(9,5).Value = 4 ; $ws1.Cells.Item(9,6).Value = 16 ; $ws1.Cells.Item(9,7).Value = 52 ; $ws1.Cells.Item(9,8).Value = 8 ; $ws1.Cells.Item(9,9).Value = 2.81 ; $ws1.Cells.Item(9,10).Value = 2.81547212600708
$ws1.Cells.Item(10,1).Value = 4 ; $ws1.Cells.Item(10,2).Value = 12 ; $ws1.Cells.Item(10,3).Value = 8 ; $ws1.Cells.Item(10,4).Value = 0 ; $ws1.Cells.Item(10,5).Value = 8 ; $ws1.Cells.Item(10,6).Value = 8 ; $ws1.Cells.Item(10,7).Value = 56.00000000000001 ; $ws1.Cells.Item(10,8).Value = 4 ; $ws1.Cells.Item(10,9).Value = 2.88 ; $ws1.Cells.Item(10,10).Value = 2.807631969451904
$ws1.Cells.Item(11,1).Value = 4 ; $ws1.Cells.Item(11,2).Value = 4 ; $ws1.Cells.Item(11,3).Value = 4 ; $ws1.Cells.Item(11,4).Value = 8 ; $ws1.Cells.Item(11,5).Value = 28 ; $ws1.Cells.Item(11,6).Value = 16 ; $ws1.Cells.Item(11,7).Value = 32 ; $ws1.Cells.Item(11,8).Value = 4 ; $ws1.Cells.Item(11,9).Value = 3.04 ; $ws1.Cells.Item(11,10).Value = 3.040929079055786
$ws1.Cells.Item(12,1).Value = 4 ; $ws1.Cells.Item(12,2).Value = 0 ; $ws1.Cells.Item(12,3).Value = 8 ; $ws1.Cells.Item(12,4).Value = 4 ; $ws1.Cells.Item(12,5).Value = 24 ; $ws1.Cells.Item(12,6).Value = 12 ; $ws1.Cells.Item(12,7).Value = 44 ; $ws1.Cells.Item(12,8).Value = 4 ; $ws1.Cells.Item(12,9).Value = 3.09 ; $ws1.Cells.Item(12,10).Value = 2.999998092651367
$ws1.Cells.Item(13,1).Value = 4 ; $ws1.Cells.Item(13,2).Value = 8 ; $ws1.Cells.Item(13,3).Value = 0 ; $ws1.Cells.Item(13,4).Value = 12 ; $ws1.Cells.Item(13,5).Value = 20 ; $ws1.Cells.Item(13,6).Value = 8 ; $ws1.Cells.Item(13,7).Value = 48 ; $ws1.Cells.Item(13,8).Value = 0 ; $ws1.Cells.Item(13,9).Value = 3.11 ; $ws1.Cells.Item(13,10).Value = 3.130991458892822
$ws1.Cells.Item(14,1).Value = 0 ; $ws1.Cells.Item(14,2).Value = 20 ; $ws1.Cells.Item(14,3).Value = 8 ; $ws1.Cells.Item(14,4).Value = 4 ; $ws1.Cells.Item(14,5).Value = 0 ; $ws1.Cells.Item(14,6).Value = 20 ; $ws1.Cells.Item(14,7).Value = 44 ; $ws1.Cells.Item(14,8).Value = 4 ; $ws1.Cells.Item(14,9).Value = 2.82 ; $ws1.Cells.Item(14,10).Value = 2.81665301322937
$ws1.Cells.Item(15,1).Value = 0 ; $ws1.Cells.Item(15,2).Value = 20 ; $ws1.Cells.Item(15,3).Value = 8 ; $ws1.Cells.Item(15,4).Value = 4 ; $ws1.Cells.Item(15,5).Value = 0 ; $ws1.Cells.Item(15,6).Value = 20 ; $ws1.Cells.Item(15,7).Value = 44 ; $ws1.Cells.Item(15,8).Value = 4 ; $ws1.Cells.Item(15,9).Value = 2.74 ; $ws1.Cells.Item(15,10).Value = 2.81665301322937
$ws1.Cells.Item(16,1).Value = 4 ; $ws1.Cells.Item(16,2).Value = 12 ; $ws1.Cells.Item(16,3).Value = 0 ; $ws1.Cells.Item(16,4).Value = 0 ; $ws1.Cells.Item(16,5).Value = 16 ; $ws1.Cells.Item(16,6).Value = 8 ; $ws1.Cells.Item(16,7).Value = 52 ; $ws1.Cells.Item(16,8).Value = 8 ; $ws1.Cells.Item(16,9).Value = 2.62 ; $ws1.Cells.Item(16,10).Value = 2.674572944641113
$ws1.Cells.Item(17,1).Value = 4 ; $ws1.Cells.Item(17,2).Value = 12 ; $ws1.Cells.Item(17,3).Value = 8 ; $ws1.Cells.Item(17,4).Value = 4 ; $ws1.Cells.Item(17,5).Value = 16 ; $ws1.Cells.Item(17,6).Value = 8 ; $ws1.Cells.Item(17,7).Value = 48 ; $ws1.Cells.Item(17,8).Value = 0 ; $ws1.Cells.Item(17,9).Value = 3.19 ; $ws1.Cells.Item(17,10).Value = 2.985782861709595
$ws1.Cells.Item(18,1).Value = 0 ; $ws1.Cells.Item(18,2).Value = 20 ; $ws1.Cells.Item(18,3).Value = 4 ; $ws1.Cells.Item(18,4).Value = 0 ; $ws1.Cells.Item(18,5).Value = 4 ; $ws1.Cells.Item(18,6).Value = 4 ; $ws1.Cells.Item(18,7).Value = 48 ; $ws1.Cells.Item(18,8).Value = 20 ; $ws1.Cells.Item(18,9).Value = 2.52 ; $ws1.Cells.Item(18,10).Value = 2.385874032974243
$ws1.Cells.Item(19,1).Value = 4 ; $ws1.Cells.Item(19,2).Value = 8 ; $ws1.Cells.Item(19,3).Value = 4 ; $ws1.Cells.Item(19,4).Value = 4 ; $ws1.Cells.Item(19,5).Value = 16 ; $ws1.Cells.Item(19,6).Value = 12 ; $ws1.Cells.Item(19,7).Value = 48 ; $ws1.Cells.Item(19,8).Value = 4 ; $ws1.Cells.Item(19,9).Value = 2.94 ; $ws1.Cells.Item(19,10).Value = 2.885764360427856
$ws1.Cells.Item(20,1).Value = 4 ; $ws1.Cells.Item(20,2).Value = 12 ; $ws1.Cells.Item(20,3).Value = 0 ; $ws1.Cells.Item(20,4).Value = 12 ; $ws1.Cells.Item(20,5).Value = 20 ; $ws1.Cells.Item(20,6).Value = 8 ; $ws1.Cells.Item(20,7).Value = 44 ; $ws1.Cells.Item(20,8).Value = 0 ; $ws1.Cells.Item(20,9).Value = 3.29 ; $ws1.Cells.Item(20,10).Value = 3.176260709762573
$ws1.Cells.Item(21,1).Value = 4 ; $ws1.Cells.Item(21,2).Value = 12 ; $ws1.Cells.Item(21,3).Value = 0 ; $ws1.Cells.Item(21,4).Value = 12 ; $ws1.Cells.Item(21,5).Value = 20 ; $ws1.Cells.Item(21,6).Value = 8 ; $ws1.Cells.Item(21,7).Value = 44 ; $ws1.Cells.Item(21,8).Value = 0 ; $ws1.Cells.Item(21,9).Value = 3.25 ; $ws1.Cells.Item(21,10).Value = 3.176260709762573
$ws1.Cells.Item(22,1).Value = 0 ; $ws1.Cells.Item(22,2).Value = 12 ; $ws1.Cells.Item(22,3).Value = 8 ; $ws1.Cells.Item(22,4).Value = 4 ; $ws1.Cells.Item(22,5).Value = 20 ; $ws1.Cells.Item(22,6).Value = 20 ; $ws1.Cells.Item(22,7).Value = 36 ; $ws1.Cells.Item(22,8).Value = 4 ; $ws1.Cells.Item(22,9).Value = 3.01 ; $ws1.Cells.Item(22,10).Value = 2.898385524749756
$ws1.Cells.Item(23,1).Value = 4 ; $ws1.Cells.Item(23,2).Value = 12 ; $ws1.Cells.Item(23,3).Value = 4 ; $ws1.Cells.Item(23,4).Value = 4 ; $ws1.Cells.Item(23,5).Value = 20 ; $ws1.Cells.Item(23,6).Value = 16 ; $ws1.Cells.Item(23,7).Value = 32 ; $ws1.Cells.Item(23,8).Value = 8 ; $ws1.Cells.Item(23,9).Value = 2.96 ; $ws1.Cells.Item(23,10).Value = 2.865798711776733
$ws1.Cells.Item(24,1).Value = 4 ; $ws1.Cells.Item(24,2).Value = 8 ; $ws1.Cells.Item(24,3).Value = 0 ; $ws1.Cells.Item(24,4).Value = 0 ; $ws1.Cells.Item(24,5).Value = 4 ; $ws1.Cells.Item(24,6).Value = 8 ; $ws1.Cells.Item(24,7).Value = 52 ; $ws1.Cells.Item(24,8).Value = 24 ; $ws1.Cells.Item(24,9).Value = 2.44 ; $ws1.Cells.Item(24,10).Value = 2.385229587554932
$ws1.Cells.Item(25,1).Value = 4 ; $ws1.Cells.Item(25,2).Value = 0 ; $ws1.Cells.Item(25,3).Value = 4 ; $ws1.Cells.Item(25,4).Value = 0 ; $ws1.Cells.Item(25,5).Value = 20 ; $ws1.Cells.Item(25,6).Value = 8 ; $ws1.Cells.Item(25,7).Value = 52 ; $ws1.Cells.Item(25,8).Value = 12 ; $ws1.Cells.Item(25,9).Value = 2.92 ; $ws1.Cells.Item(25,10).Value = 2.732897520065308
$ws1.Cells.Item(26,1).Value = 4 ; $ws1.Cells.Item(26,2).Value = 8 ; $ws1.Cells.Item(26,3).Value = 12 ; $ws1.Cells.Item(26,4).Value = 4 ; $ws1.Cells.Item(26,5).Value = 4 ; $ws1.Cells.Item(26,6).Value = 24 ; $ws1.Cells.Item(26,7).Value = 39.99999999999999 ; $ws1.Cells.Item(26,8).Value = 4 ; $ws1.Cells.Item(26,9).Value = 2.9 ; $ws1.Cells.Item(26,10).Value = 2.886198997497559
$ws1.Cells.Item(27,1).Value = 4 ; $ws1.Cells.Item(27,2).Value = 0 ; $ws1.Cells.Item(27,3).Value = 12 ; $ws1.Cells.Item(27,4).Value = 16 ; $ws1.Cells.Item(27,5).Value = 4 ; $ws1.Cells.Item(27,6).Value = 12 ; $ws1.Cells.Item(27,7).Value = 52 ; $ws1.Cells.Item(27,8).Value = 0 ; $ws1.Cells.Item(27,9).Value = 2.96 ; $ws1.Cells.Item(27,10).Value = 2.93699836730957
$ws1.Cells.Item(28,1).Value = 0 ; $ws1.Cells.Item(28,2).Value = 16 ; $ws1.Cells.Item(28,3).Value = 8 ; $ws1.Cells.Item(28,4).Value = 4 ; $ws1.Cells.Item(28,5).Value = 16 ; $ws1.Cells.Item(28,6).Value = 20 ; $ws1.Cells.Item(28,7).Value = 36 ; $ws1.Cells.Item(28,8).Value = 0 ; $ws1.Cells.Item(28,9).Value = 3.02 ; $ws1.Cells.Item(28,10).Value = 2.911064624786377
$ws1.Cells.Item(29,1).Value = 4 ; $ws1.Cells.Item(29,2).Value = 12 ; $ws1.Cells.Item(29,3).Value = 8 ; $ws1.Cells.Item(29,4).Value = 0 ; $ws1.Cells.Item(29,5).Value = 8 ; $ws1.Cells.Item(29,6).Value = 8 ; $ws1.Cells.Item(29,7).Value = 56.00000000000001 ; $ws1.Cells.Item(29,8).Value = 4 ; $ws1.Cells.Item(29,9).Value = 2.79 ; $ws1.Cells.Item(29,10).Value = 2.807631969451904
$ws1.Cells.Item(30,1).Value = 4 ; $ws1.Cells.Item(30,2).Value = 0 ; $ws1.Cells.Item(30,3).Value = 4 ; $ws1.Cells.Item(30,4).Value = 4 ; $ws1.Cells.Item(30,5).Value = 16 ; $ws1.Cells.Item(30,6).Value = 12 ; $ws1.Cells.Item(30,7).Value = 56.00000000000001 ; $ws1.Cells.Item(30,8).Value = 4 ; $ws1.Cells.Item(30,9).Value = 2.85 ; $ws1.Cells.Item(30,10).Value = 2.898507595062256
$ws1.Cells.Item(31,1).Value = 0 ; $ws1.Cells.Item(31,2).Value = 4 ; $ws1.Cells.Item(31,3).Value = 8 ; $ws1.Cells.Item(31,4).Value = 4 ; $ws1.Cells.Item(31,5).Value = 16 ; $ws1.Cells.Item(31,6).Value = 20 ; $ws1.Cells.Item(31,7).Value = 44 ; $ws1.Cells.Item(31,8).Value = 4 ; $ws1.Cells.Item(31,9).Value = 2.94 ; $ws1.Cells.Item(31,10).Value = 2.920213460922241
$ws1.Cells.Item(32,1).Value = 0 ; $ws1.Cells.Item(32,2).Value = 12 ; $ws1.Cells.Item(32,3).Value = 4 ; $ws1.Cells.Item(32,4).Value = 0 ; $ws1.Cells.Item(32,5).Value = 4 ; $ws1.Cells.Item(32,6).Value = 16 ; $ws1.Cells.Item(32,7).Value = 39.99999999999999 ; $ws1.Cells.Item(32,8).Value = 24 ; $ws1.Cells.Item(32,9).Value = 2.51 ; $ws1.Cells.Item(32,10).Value = 2.413007974624634
$ws1.Cells.Item(33,1).Value = 0 ; $ws1.Cells.Item(33,2).Value = 0 ; $ws1.Cells.Item(33,3).Value = 8 ; $ws1.Cells.Item(33,4).Value = 4 ; $ws1.Cells.Item(33,5).Value = 16 ; $ws1.Cells.Item(33,6).Value = 12 ; $ws1.Cells.Item(33,7).Value = 52 ; $ws1.Cells.Item(33,8).Value = 8 ; $ws1.Cells.Item(33,9).Value = 2.99 ; $ws1.Cells.Item(33,10).Value = 2.89897632598877
$ws1.Cells.Item(34,1).Value = 0 ; $ws1.Cells.Item(34,2).Value = 16 ; $ws1.Cells.Item(34,3).Value = 0 ; $ws1.Cells.Item(34,4).Value = 4 ; $ws1.Cells.Item(34,5).Value = 20 ; $ws1.Cells.Item(34,6).Value = 20 ; $ws1.Cells.Item(34,7).Value = 28 ; $ws1.Cells.Item(34,8).Value = 12 ; $ws1.Cells.Item(34,9).Value = 3.38 ; $ws1.Cells.Item(34,10).Value = 2.86850380897522
$ws1.Cells.Item(35,1).Value = 4 ; $ws1.Cells.Item(35,2).Value = 16 ; $ws1.Cells.Item(35,3).Value = 8 ; $ws1.Cells.Item(35,4).Value = 0 ; $ws1.Cells.Item(35,5).Value = 12 ; $ws1.Cells.Item(35,6).Value = 0 ; $ws1.Cells.Item(35,7).Value = 48 ; $ws1.Cells.Item(35,8).Value = 12 ; $ws1.Cells.Item(35,9).Value = 2.56 ; $ws1.Cells.Item(35,10).Value = 2.477823972702026
$ws1.Cells.Item(36,1).Value = 0 ; $ws1.Cells.Item(36,2).Value = 0 ; $ws1.Cells.Item(36,3).Value = 8 ; $ws1.Cells.Item(36,4).Value = 4 ; $ws1.Cells.Item(36,5).Value = 16 ; $ws1.Cells.Item(36,6).Value = 12 ; $ws1.Cells.Item(36,7).Value = 52 ; $ws1.Cells.Item(36,8).Value = 8 ; $ws1.Cells.Item(36,9).Value = 2.82 ; $ws1.Cells.Item(36,10).Value = 2.89897632598877
$ws1.Cells.Item(37,1).Value = 4 ; $ws1.Cells.Item(37,2).Value = 12 ; $ws1.Cells.Item(37,3).Value = 4 ; $ws1.Cells.Item(37,4).Value = 4 ; $ws1.Cells.Item(37,5).Value = 12 ; $ws1.Cells.Item(37,6).Value = 16 ; $ws1.Cells.Item(37,7).Value = 36 ; $ws1.Cells.Item(37,8).Value = 12 ; $ws1.Cells.Item(37,9).Value = 2.86 ; $ws1.Cells.Item(37,10).Value = 2.729341506958008
$ws1.Cells.Item(38,1).Value = 0 ; $ws1.Cells.Item(38,2).Value = 24 ; $ws1.Cells.Item(38,3).Value = 8 ; $ws1.Cells.Item(38,4).Value = 0 ; $ws1.Cells.Item(38,5).Value = 12 ; $ws1.Cells.Item(38,6).Value = 16 ; $ws1.Cells.Item(38,7).Value = 32 ; $ws1.Cells.Item(38,8).Value = 8 ; $ws1.Cells.Item(38,9).Value = 2.93 ; $ws1.Cells.Item(38,10).Value = 2.715706348419189
$ws1.Cells.Item(39,1).Value = 4 ; $ws1.Cells.Item(39,2).Value = 32 ; $ws1.Cells.Item(39,3).Value = 8 ; $ws1.Cells.Item(39,4).Value = 4 ; $ws1.Cells.Item(39,5).Value = 4 ; $ws1.Cells.Item(39,6).Value = 20 ; $ws1.Cells.Item(39,7).Value = 24 ; $ws1.Cells.Item(39,8).Value = 4 ; $ws1.Cells.Item(39,9).Value = 2.84 ; $ws1.Cells.Item(39,10).Value = 2.806526660919189
$ws1.Cells.Item(40,1).Value = 0 ; $ws1.Cells.Item(40,2).Value = 24 ; $ws1.Cells.Item(40,3).Value = 8 ; $ws1.Cells.Item(40,4).Value = 0 ; $ws1.Cells.Item(40,5).Value = 12 ; $ws1.Cells.Item(40,6).Value = 16 ; $ws1.Cells.Item(40,7).Value = 32 ; $ws1.Cells.Item(40,8).Value = 8 ; $ws1.Cells.Item(40,9).Value = 2.94 ; $ws1.Cells.Item(40,10).Value = 2.715706348419189
$ws1.Cells.Item(41,1).Value = 4 ; $ws1.Cells.Item(41,2).Value = 40 ; $ws1.Cells.Item(41,3).Value = 0 ; $ws1.Cells.Item(41,4).Value = 0 ; $ws1.Cells.Item(41,5).Value = 12 ; $ws1.Cells.Item(41,6).Value = 4 ; $ws1.Cells.Item(41,7).Value = 36 ; $ws1.Cells.Item(41,8).Value = 4 ; $ws1.Cells.Item(41,9).Value = 3.16 ; $ws1.Cells.Item(41,10).Value = 3.167737245559692
$ws1.Cells.Item(42,1).Value = 0 ; $ws1.Cells.Item(42,2).Value = 16 ; $ws1.Cells.Item(42,3).Value = 0 ; $ws1.Cells.Item(42,4).Value = 4 ; $ws1.Cells.Item(42,5).Value = 20 ; $ws1.Cells.Item(42,6).Value = 20 ; $ws1.Cells.Item(42,7).Value = 28 ; $ws1.Cells.Item(42,8).Value = 12 ; $ws1.Cells.Item(42,9).Value = 2.72 ; $ws1.Cells.Item(42,10).Value = 2.86850380897522
$ws1.Cells.Item(43,1).Value = 4 ; $ws1.Cells.Item(43,2).Value = 20 ; $ws1.Cells.Item(43,3).Value = 4 ; $ws1.Cells.Item(43,4).Value = 4 ; $ws1.Cells.Item(43,5).Value = 16 ; $ws1.Cells.Item(43,6).Value = 0 ; $ws1.Cells.Item(43,7).Value = 52 ; $ws1.Cells.Item(43,8).Value = 0 ; $ws1.Cells.Item(43,9).Value = 3.21 ; $ws1.Cells.Item(43,10).Value = 3.116974115371704
$ws1.Cells.Item(44,1).Value = 4 ; $ws1.Cells.Item(44,2).Value = 12 ; $ws1.Cells.Item(44,3).Value = 0 ; $ws1.Cells.Item(44,4).Value = 0 ; $ws1.Cells.Item(44,5).Value = 16 ; $ws1.Cells.Item(44,6).Value = 8 ; $ws1.Cells.Item(44,7).Value = 52 ; $ws1.Cells.Item(44,8).Value = 8 ; $ws1.Cells.Item(44,9).Value = 2.56 ; $ws1.Cells.Item(44,10).Value = 2.674572944641113

# Sheet2 (Test Results) updates: rows 2-11, columns A-J
$ws2.Cells.Item(2,1).Value = 4 ; $ws2.Cells.Item(2,2).Value = 0 ; $ws2.Cells.Item(2,3).Value = 12 ; $ws2.Cells.Item(2,4).Value = 16 ; $ws2.Cells.Item(2,5).Value = 4 ; $ws2.Cells.Item(2,6).Value = 12 ; $ws2.Cells.Item(2,7).Value = 52 ; $ws2.Cells.Item(2,8).Value = 0 ; $ws2.Cells.Item(2,9).Value = 2.96 ; $ws2.Cells.Item(2,10).Value = 2.93699836730957
$ws2.Cells.Item(3,1).Value = 4 ; $ws2.Cells.Item(3,2).Value = 40 ; $ws2.Cells.Item(3,3).Value = 0 ; $ws2.Cells.Item(3,4).Value = 0 ; $ws2.Cells.Item(3,5).Value = 12 ; $ws2.Cells.Item(3,6).Value = 4 ; $ws2.Cells.Item(3,7).Value = 36 ; $ws2.Cells.Item(3,8).Value = 4 ; $ws2.Cells.Item(3,9).Value = 3.16 ; $ws2.Cells.Item(3,10).Value = 3.167737245559692
$ws2.Cells.Item(4,1).Value = 4 ; $ws2.Cells.Item(4,2).Value = 12 ; $ws2.Cells.Item(4,3).Value = 4 ; $ws2.Cells.Item(4,4).Value = 4 ; $ws2.Cells.Item(4,5).Value = 12 ; $ws2.Cells.Item(4,6).Value = 16 ; $ws2.Cells.Item(4,7).Value = 36 ; $ws2.Cells.Item(4,8).Value = 12 ; $ws2.Cells.Item(4,9).Value = 2.86 ; $ws2.Cells.Item(4,10).Value = 2.729341506958008
$ws2.Cells.Item(5,1).Value = 0 ; $ws2.Cells.Item(5,2).Value = 20 ; $ws2.Cells.Item(5,3).Value = 8 ; $ws2.Cells.Item(5,4).Value = 4 ; $ws2.Cells.Item(5,5).Value = 0 ; $ws2.Cells.Item(5,6).Value = 20 ; $ws2.Cells.Item(5,7).Value = 44 ; $ws2.Cells.Item(5,8).Value = 4 ; $ws2.Cells.Item(5,9).Value = 2.82 ; $ws2.Cells.Item(5,10).Value = 2.81665301322937
$ws2.Cells.Item(6,1).Value = 4 ; $ws2.Cells.Item(6,2).Value = 12 ; $ws2.Cells.Item(6,3).Value = 8 ; $ws2.Cells.Item(6,4).Value = 4 ; $ws2.Cells.Item(6,5).Value = 16 ; $ws2.Cells.Item(6,6).Value = 8 ; $ws2.Cells.Item(6,7).Value = 48 ; $ws2.Cells.Item(6,8).Value = 0 ; $ws2.Cells.Item(6,9).Value = 3.19 ; $ws2.Cells.Item(6,10).Value = 2.985782861709595
$ws2.Cells.Item(7,1).Value = 4 ; $ws2.Cells.Item(7,2).Value = 0 ; $ws2.Cells.Item(7,3).Value = 8 ; $ws2.Cells.Item(7,4).Value = 4 ; $ws2.Cells.Item(7,5).Value = 24 ; $ws2.Cells.Item(7,6).Value = 12 ; $ws2.Cells.Item(7,7).Value = 44 ; $ws2.Cells.Item(7,8).Value = 4 ; $ws2.Cells.Item(7,9).Value = 3.09 ; $ws2.Cells.Item(7,10).Value = 2.999998092651367
$ws2.Cells.Item(8,1).Value = 0 ; $ws2.Cells.Item(8,2).Value = 4 ; $ws2.Cells.Item(8,3).Value = 8 ; $ws2.Cells.Item(8,4).Value = 4 ; $ws2.Cells.Item(8,5).Value = 16 ; $ws2.Cells.Item(8,6).Value = 20 ; $ws2.Cells.Item(8,7).Value = 44 ; $ws2.Cells.Item(8,8).Value = 4 ; $ws2.Cells.Item(8,9).Value = 2.94 ; $ws2.Cells.Item(8,10).Value = 2.920213460922241
$ws2.Cells.Item(9,1).Value = 4 ; $ws2.Cells.Item(9,2).Value = 12 ; $ws2.Cells.Item(9,3).Value = 0 ; $ws2.Cells.Item(9,4).Value = 12 ; $ws2.Cells.Item(9,5).Value = 20 ; $ws2.Cells.Item(9,6).Value = 8 ; $ws2.Cells.Item(9,7).Value = 44 ; $ws2.Cells.Item(9,8).Value = 0 ; $ws2.Cells.Item(9,9).Value = 3.25 ; $ws2.Cells.Item(9,10).Value = 3.176260709762573
$ws2.Cells.Item(10,1).Value = 0 ; $ws2.Cells.Item(10,2).Value = 20 ; $ws2.Cells.Item(10,3).Value = 8 ; $ws2.Cells.Item(10,4).Value = 4 ; $ws2.Cells.Item(10,5).Value = 0 ; $ws2.Cells.Item(10,6).Value = 20 ; $ws2.Cells.Item(10,7).Value = 44 ; $ws2.Cells.Item(10,8).Value = 4 ; $ws2.Cells.Item(10,9).Value = 2.74 ; $ws2.Cells.Item(10,10).Value = 2.81665301322937
$ws2.Cells.Item(11,1).Value = 4 ; $ws2.Cells.Item(11,2).Value = 8 ; $ws2.Cells.Item(11,3).Value = 4 ; $ws2.Cells.Item(11,4).Value = 4 ; $ws2.Cells.Item(11,5).Value = 16 ; $ws2.Cells.Item(11,6).Value = 12 ; $ws2.Cells.Item(11,7).Value = 48 ; $ws2.Cells.Item(11,8).Value = 4 ; $ws2.Cells.Item(11,9).Value = 2.94 ; $ws2.Cells.Item(11,10).Value = 2.885764360427856
